$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, pushing existing rows 103-153 down to 104-154.
$ws.Rows(103).Insert()

# Populate the newly inserted row 103 with a new data record
# (same Mercado/Region/Categoria/Calidad/Codreg/CategoriaID/Clasificacion/KgUnidades
# as the surrounding rows for this subset, with new Fecha/Variedad/Volumen/Precios/
# Unidad/Origen values).
$ws.Range("A103").Value2 = 4
$ws.Range("B103").Value2 = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C103").Value2 = 'Los Lagos'
$ws.Range("D103").Value2 = 45062
$ws.Range("E103").Value2 = 10
$ws.Range("F103").Value2 = 100112022
$ws.Range("G103").Value2 = 'Arveja Verde'
$ws.Range("H103").Value2 = 'Perfection'
$ws.Range("I103").Value2 = 'Primera'
$ws.Range("J103").Value2 = 35
$ws.Range("K103").Value2 = 43000
$ws.Range("L103").Value2 = 43000
$ws.Range("M103").Value2 = 43000
$ws.Range("N103").Value2 = '$/malla 25 kilos'
$ws.Range("O103").Value2 = 'Provincia de Huasco'
$ws.Range("P103").Value2 = 1720
$ws.Range("Q103").Value2 = 25
$ws.Range("R103").Value2 = 'Hortaliza'
